# The doc has two tag/id blocks that were previously split into three
# differently-formatted runs each:
#   "<id>"  (Courier New, 7f6000)  +  "p033v_1"  (plain/black)  +  "</id>" (Courier New, 7f6000)
# The freshly re-downloaded source now represents each id block as a single
# run: "<id>p033v_1</id>" (and likewise for p033v_2), carrying the Courier
# New / 7f6000 formatting of the original opening-tag run.
#
# Using Find/Replace across the whole (already-contiguous) text collapses
# those three runs into one run, taking on the formatting of the first
# character of the match - exactly the opening "<id>" run's formatting.

$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()

$find.Execute("<id>p033v_1</id>", $true, $true, $false, $false, $false, $true, 1, $false, "<id>p033v_1</id>", 2)
$find.Execute("<id>p033v_2</id>", $true, $true, $false, $false, $false, $true, 1, $false, "<id>p033v_2</id>", 2)
